$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.757.24'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '2.081.13'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.650'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '53.56'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.83'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.365'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0760'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.94'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.885'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('D15').Value = '2.384.13'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('E16').Value = '  -3.75%  '
$ws.Range('D17').Value = '2.047.04'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '36.732.73'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('D21').Value = '0.0₃0877'
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.60'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.98%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.123'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('E32').Value = '  +5.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0605'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('E35').Value = '  +5.91%  '
$ws.Range('E37').Value = '  +4.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0830'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.30%  '
$ws.Range('E39').Value = '  -5.04%  '
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.84'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0947'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('E45').Value = '  -12.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.30%  '
$ws.Range('D47').Value = '1.357.76'
$ws.Range('E47').Value = '  +6.46%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.28'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.95%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('D51').Value = '2.267.57'
$ws.Range('E51').Value = '  +1.62%  '
